$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers but must stay as literal text
# (matching the source data which mixes "61.741.36"-style and plain decimal strings).
# Force text format first so Excel does not convert them to numeric values.
$textCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D12", "D13", "D16", "D21", "D22", "D23", "D24", "D26", "D27", "D28", "D29", "D30", "D33", "D34", "D35", "D36", "D37", "D39", "D40", "D41", "D42", "D45", "D46", "D47", "D48", "D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "61.813.55"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").Value = "3.418.89"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").Value = "409.75"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").Value = "128.99"
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("D7").Value = "0.630"
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "0.729"
$ws.Range("E9").Value = "  -3.56%  "
$ws.Range("D10").Value = "0.139"
$ws.Range("E10").Value = "  -1.00%  "
$ws.Range("D11").Value = "43.32"
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("D12").Value = "0.0000222"
$ws.Range("E12").Value = "  +14.41%  "
$ws.Range("D13").Value = "9.25"
$ws.Range("E13").Value = "  +5.24%  "
$ws.Range("D14").Value = "3.958.39"
$ws.Range("E14").Value = "  -0.70%  "
$ws.Range("E15").Value = "  +0.30%  "
$ws.Range("D16").Value = "21.13"
$ws.Range("E16").Value = "  +4.06%  "
$ws.Range("D17").Value = "3.413.48"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("E18").Value = "  +7.87%  "
$ws.Range("E19").Value = "  +3.15%  "
$ws.Range("D20").Value = "61.753.30"
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("D21").Value = "485.19"
$ws.Range("E21").Value = "  +29.88%  "
$ws.Range("D22").Value = "91.46"
$ws.Range("E22").Value = "  +4.46%  "
$ws.Range("D23").Value = "3.30"
$ws.Range("E23").Value = "  +3.72%  "
$ws.Range("D24").Value = "13.51"
$ws.Range("E24").Value = "  +1.67%  "
$ws.Range("D26").Value = "34.50"
$ws.Range("E26").Value = "  +8.80%  "
$ws.Range("D27").Value = "9.25"
$ws.Range("E27").Value = "  +9.26%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "7.61"
$ws.Range("E28").Value = "  -1.14%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "12.11"
$ws.Range("E29").Value = "  +2.37%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "2.69"
$ws.Range("E30").Value = "  -1.42%  "
$ws.Range("E31").Value = "  -1.47%  "
$ws.Range("E32").Value = "  -2.07%  "
$ws.Range("D33").Value = "41.92"
$ws.Range("E33").Value = "  -4.72%  "
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").Value = "58.70"
$ws.Range("E35").Value = "  +12.50%  "
$ws.Range("D36").Value = "0.0497"
$ws.Range("E36").Value = "  +0.81%  "
$ws.Range("D37").Value = "0.998"
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("E38").Value = "  +2.94%  "
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").Value = "0.136"
$ws.Range("E39").Value = "  +3.22%  "
$ws.Range("B40").Value = "WEMIXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").Value = "2.74"
$ws.Range("E40").Value = "  +17.38%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "146.07"
$ws.Range("E41").Value = "  +1.88%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "2.92"
$ws.Range("E42").Value = "  +0.49%  "
$ws.Range("E43").Value = "  +1.88%  "
$ws.Range("E44").Value = "  +5.24%  "
$ws.Range("D45").Value = "4.35"
$ws.Range("E45").Value = "  +8.62%  "
$ws.Range("B46").Value = "Celestia"
$ws.Range("C46").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D46").Value = "16.69"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").Value = "2.33"
$ws.Range("E47").Value = "  +20.39%  "
$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D48").Value = "118.77"
$ws.Range("E48").Value = "  +27.79%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "22.73"
$ws.Range("E49").Value = "  +3.88%  "
$ws.Range("E50").Value = "  +16.98%  "
$ws.Range("D51").Value = "2.133.63"
$ws.Range("E51").Value = "  +1.00%  "

# Restore default (General) formatting on those cells so no stray number format
# sticks around on what should remain plain, unstyled text cells.
foreach ($addr in $textCells) {
    $ws.Range($addr).ClearFormats()
}
